$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing a Text-formatted
# result (so digit/decimal-looking strings like "253.74" are not
# auto-coerced into numbers by Excel), then drop back to the default
# "Normal" style so no stray number-format is left on the cell.
function Set-TextValue {
    param($CellRef, $Val)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "37.169.44"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.073.42"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "253.74"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("E6").Value = "  +2.20%  "
Set-TextValue "D7" "59.29"
$ws.Range("E7").Value = "  +12.61%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +4.13%  "
Set-TextValue "D10" "61.57"
$ws.Range("E10").Value = "  -0.16%  "
Set-TextValue "D11" "0.0794"
$ws.Range("E11").Value = "  +6.88%  "
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("E13").Value = "  +6.87%  "
$ws.Range("D14").Value = "2.379.20"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("E15").Value = "  -1.70%  "
Set-TextValue "D16" "5.56"
$ws.Range("E16").Value = "  +8.17%  "
$ws.Range("D17").Value = "2.072.83"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "37.140.37"
$ws.Range("E18").Value = "  -0.34%  "
Set-TextValue "D19" "15.78"
$ws.Range("E19").Value = "  +11.07%  "
Set-TextValue "D20" "74.85"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "0.0₃0934"
$ws.Range("E21").Value = "  +11.13%  "
Set-TextValue "D22" "5.47"
$ws.Range("E22").Value = "  +4.41%  "
Set-TextValue "D23" "240.48"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("E26").Value = "  +14.66%  "
Set-TextValue "D27" "170.14"
$ws.Range("E27").Value = "  -0.29%  "
Set-TextValue "D28" "9.37"
$ws.Range("E28").Value = "  +1.53%  "
Set-TextValue "D29" "20.35"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("E31").Value = "  +7.26%  "
Set-TextValue "D32" "4.81"
$ws.Range("E32").Value = "  +6.88%  "
Set-TextValue "D33" "0.0636"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("E34").Value = "  +9.00%  "
Set-TextValue "D35" "0.0915"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +1.00%  "
Set-TextValue "D38" "0.119"
$ws.Range("E38").Value = "  +29.29%  "
Set-TextValue "D39" "1.77"
$ws.Range("E39").Value = "  -4.74%  "
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0228"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D42" "17.91"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  +0.51%  "
Set-TextValue "D44" "98.84"
$ws.Range("E44").Value = "  +0.03%  "
Set-TextValue "D45" "4.38"
$ws.Range("E45").Value = "  +12.69%  "
Set-TextValue "D46" "2.82"
$ws.Range("E46").Value = "  +2.58%  "
Set-TextValue "D47" "4.52"
$ws.Range("E47").Value = "  +13.00%  "
Set-TextValue "D48" "2.48"
$ws.Range("E48").Value = "  +8.32%  "
$ws.Range("D49").Value = "1.306.09"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -0.88%  "
